$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.303.58"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "3.830.92"
$ws.Range("E3").Value = "  +4.01%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'414.86"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "'134.46"
$ws.Range("E6").Value = "  +3.27%  "
$ws.Range("D7").Value = "3.820.34"
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("D8").Value = "'0.619"
$ws.Range("E8").Value = "  -3.11%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'0.748"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "'0.172"
$ws.Range("E11").Value = "  -3.75%  "
$ws.Range("D12").Value = "'0.0000381"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "'41.46"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").Value = "'10.24"
$ws.Range("E14").Value = "  -3.71%  "
$ws.Range("D15").Value = "4.382.13"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "'14.90"
$ws.Range("E16").Value = "  +16.03%  "
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "3.808.14"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").Value = "'19.69"
$ws.Range("E19").Value = "  -3.71%  "
$ws.Range("D20").Value = "67.262.09"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").Value = "'1.08"
$ws.Range("E21").Value = "  -2.99%  "
$ws.Range("D22").Value = "'418.69"
$ws.Range("E22").Value = "  -4.36%  "
$ws.Range("D23").Value = "'14.90"
$ws.Range("E23").Value = "  -9.10%  "
$ws.Range("D24").Value = "'86.76"
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("D25").Value = "'3.10"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'36.65"
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'5.71"
$ws.Range("E27").Value = "  +14.69%  "
$ws.Range("D28").Value = "'3.17"
$ws.Range("E28").Value = "  -3.45%  "
$ws.Range("D29").Value = "'9.50"
$ws.Range("E29").Value = "  -7.08%  "
$ws.Range("D30").Value = "'693.44"
$ws.Range("E30").Value = "  +6.87%  "
$ws.Range("D31").Value = "'0.123"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").Value = "'12.52"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("D33").Value = "'2.70"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").Value = "'7.39"
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("D35").Value = "'0.156"
$ws.Range("E35").Value = "  -4.78%  "
$ws.Range("D36").Value = "'39.55"
$ws.Range("E36").Value = "  -4.60%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "'55.79"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("D39").Value = "0.0₃0773"
$ws.Range("E39").Value = "  +3.70%  "
$ws.Range("D40").Value = "'0.0466"
$ws.Range("E40").Value = "  -4.98%  "
$ws.Range("D41").Value = "'2.92"
$ws.Range("E41").Value = "  -8.44%  "
$ws.Range("D42").Value = "'0.994"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'0.138"
$ws.Range("E43").Value = "  -5.91%  "
$ws.Range("D44").Value = "'27.74"
$ws.Range("E44").Value = "  -13.27%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'147.26"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.18"
$ws.Range("E46").Value = "  +20.02%  "
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").Value = "'3.33"
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'2.11"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.90"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'4.35"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'2.59"
$ws.Range("E51").Value = "  +0.68%  "
